$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# K1: seas_id_x -> seas_id
$ws.Range("K1").Value = "seas_id"
# L1: player_id -> player_id_x
$ws.Range("L1").Value = "player_id_x"
# M1: season_x (unchanged)
# N1: seas_id_y -> season_ending_year_y
$ws.Range("N1").Value = "season_ending_year_y"
# O1: season_ending_year_y -> player_id_y
$ws.Range("O1").Value = "player_id_y"
# P1: player_y (unchanged)

# --- Update data rows 2-43 ---
# Column N becomes the season_ending_year_y text value (what used to be in column O)
# Column O becomes a brand new numeric player_id_y value

$seasonEndingYear = @{
    2 = "2024"; 3 = "2023"; 4 = "2022"; 5 = "2021"; 6 = "2020"; 7 = "2019";
    8 = "2018"; 9 = "2017"; 10 = "2016"; 11 = "2015"; 12 = "2014"; 13 = "2013";
    14 = "2012"; 15 = "2011"; 16 = "2010"; 17 = "2009"; 18 = "2008"; 19 = "2007";
    20 = "2006"; 21 = "2005"; 22 = "2004"; 23 = "2003"; 24 = "2002"; 25 = "2001";
    26 = "2000"; 27 = "1999"; 28 = "1998"; 29 = "1997"; 30 = "1996"; 31 = "1995";
    32 = "1994"; 33 = "1993"; 34 = "1992"; 35 = "1991"; 36 = "1990"; 37 = "1989";
    38 = "1988"; 39 = "1987"; 40 = "1986"; 41 = "1985"; 42 = "1984"; 43 = "1983";
}

$playerIdY = @{
    2 = 3832; 3 = 3412; 4 = 5013; 5 = 2865; 6 = 3776; 7 = 3354;
    8 = 3354; 9 = 1677; 10 = 2332; 11 = 3354; 12 = 2332; 13 = 2220;
    14 = 2354; 15 = 3192; 16 = 2332; 17 = 2427; 18 = 3435; 19 = 3263;
    20 = 3709; 21 = 343; 22 = 216; 23 = 558; 24 = 1023; 25 = 21;
    26 = 4339; 27 = 1173; 28 = 1144; 29 = 2794; 30 = 4901; 31 = 235;
    32 = 1289; 33 = 978; 34 = 1355; 35 = 1355; 36 = 4288; 37 = 1621;
    38 = 4425; 39 = 4288; 40 = 422; 41 = 3107; 42 = 3107; 43 = 559;
}

for ($r = 2; $r -le 43; $r++) {
    # Leading apostrophe forces Excel to store the numeric-looking year as
    # text (shared string) rather than re-inferring it as a number; reset
    # the style back to Normal afterwards so no stray number-format/style
    # is left behind on the cell.
    $ws.Range("N$r").Value = "'" + $seasonEndingYear[$r]
    $ws.Range("N$r").Style = "Normal"
    $ws.Range("O$r").Value = $playerIdY[$r]
}
